$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.833.55'
$ws.Range("E2").Value = '  -0.40%  '
# Row 3
$ws.Range("D3").Value = '3.806.87'
$ws.Range("E3").Value = '  -1.44%  '
# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.09%  '
# Row 5
$ws.Range("D5").Value = '''705.03'
$ws.Range("E5").Value = '  +0.48%  '
# Row 6
$ws.Range("D6").Value = '''170.06'
$ws.Range("E6").Value = '  -2.13%  '
# Row 7
$ws.Range("D7").Value = '3.807.01'
$ws.Range("E7").Value = '  -1.34%  '
# Row 8
$ws.Range("E8").Value = '  +0.00%  '
# Row 9
$ws.Range("E9").Value = '  -0.79%  '
# Row 10
$ws.Range("D10").Value = '''0.161'
$ws.Range("E10").Value = '  -1.49%  '
# Row 11
$ws.Range("D11").Value = '''7.55'
$ws.Range("E11").Value = '  +5.43%  '
# Row 12
$ws.Range("D12").Value = '''0.458'
$ws.Range("E12").Value = '  -0.60%  '
# Row 13
$ws.Range("E13").Value = '  -2.62%  '
# Row 14
$ws.Range("D14").Value = '''35.93'
$ws.Range("E14").Value = '  -1.73%  '
# Row 15
$ws.Range("D15").Value = '4.449.02'
$ws.Range("E15").Value = '  -1.46%  '
# Row 16
$ws.Range("D16").Value = '3.794.59'
$ws.Range("E16").Value = '  -1.79%  '
# Row 17
$ws.Range("D17").Value = '70.846.75'
$ws.Range("E17").Value = '  -0.53%  '
# Row 18
$ws.Range("E18").Value = '  +0.01%  '
# Row 19
$ws.Range("D19").Value = '''7.10'
$ws.Range("E19").Value = '  -1.89%  '
# Row 20
$ws.Range("D20").Value = '''17.34'
$ws.Range("E20").Value = '  -2.24%  '
# Row 21
$ws.Range("D21").Value = '''496.33'
$ws.Range("E21").Value = '  -0.72%  '
# Row 22
$ws.Range("D22").Value = '''10.69'
$ws.Range("E22").Value = '  -4.75%  '
# Row 23
$ws.Range("E23").Value = '  +0.02%  '
# Row 24
$ws.Range("D24").Value = '''84.17'
$ws.Range("E24").Value = '  -0.99%  '
# Row 25
$ws.Range("D25").Value = '''0.0000143'
$ws.Range("E25").Value = '  -1.70%  '
# Row 26
$ws.Range("D26").Value = '3.955.35'
$ws.Range("E26").Value = '  -1.33%  '
# Row 27
$ws.Range("D27").Value = '''12.08'
$ws.Range("E27").Value = '  -1.93%  '
# Row 28
$ws.Range("D28").Value = '''10.32'
$ws.Range("E28").Value = '  -4.33%  '
# Row 29
$ws.Range("E29").Value = '  +0.04%  '
# Row 30
$ws.Range("D30").Value = '''2.02'
$ws.Range("E30").Value = '  -6.09%  '
# Row 31
$ws.Range("D31").Value = '''3.03'
$ws.Range("E31").Value = '  -6.00%  '
# Row 32
$ws.Range("E32").Value = '  -0.26%  '
# Row 33
$ws.Range("E33").Value = '  -3.97%  '
# Row 34
$ws.Range("D34").Value = '''29.05'
$ws.Range("E34").Value = '  -2.30%  '
# Row 35
$ws.Range("D35").Value = '''0.176'
$ws.Range("E35").Value = '  -2.77%  '
# Row 36
$ws.Range("D36").Value = '3.775.67'
$ws.Range("E36").Value = '  -1.04%  '
# Row 37
$ws.Range("E37").Value = '  -0.11%  '
# Row 38
$ws.Range("D38").Value = '''9.08'
$ws.Range("E38").Value = '  -2.19%  '
# Row 39
$ws.Range("E39").Value = '  -3.83%  '
# Row 40
$ws.Range("D40").Value = '''2.37'
$ws.Range("E40").Value = '  -3.05%  '
# Row 41
$ws.Range("D41").Value = '''1.02'
$ws.Range("E41").Value = '  -1.85%  '
# Row 42
$ws.Range("D42").Value = '''5.94'
$ws.Range("E42").Value = '  -1.66%  '
# Row 43
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  -0.01%  '
# Row 44
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '''3.23'
$ws.Range("E44").Value = '  -5.10%  '
# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''1.00'
$ws.Range("E45").Value = '  -0.01%  '
# Row 46
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").Value = '''0.000319'
$ws.Range("E46").Value = '  +0.12%  '
# Row 47
$ws.Range("D47").Value = '''166.57'
$ws.Range("E47").Value = '  +1.70%  '
# Row 48
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").Value = '''49.01'
$ws.Range("E48").Value = '  +0.42%  '
# Row 49
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").Value = '''425.99'
$ws.Range("E49").Value = '  +2.02%  '
# Row 50
$ws.Range("E50").Value = '  -0.85%  '
# Row 51
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = '''0.293'
$ws.Range("E51").Value = '  -3.47%  '
